# Update the "想去人数" (want-to-go count) figures that changed between
# the previous data pull and the latest one (456a3b4), on both the
# "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 759
$wsExhibit.Range("F10").Value = 104
$wsExhibit.Range("F12").Value = 98
$wsExhibit.Range("F13").Value = 302
$wsExhibit.Range("F14").Value = 413
$wsExhibit.Range("F15").Value = 492
$wsExhibit.Range("F17").Value = 11146
$wsExhibit.Range("F18").Value = 5328

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 759
$wsAll.Range("F12").Value = 104
$wsAll.Range("F14").Value = 98
$wsAll.Range("F15").Value = 302
$wsAll.Range("F16").Value = 413
$wsAll.Range("F17").Value = 492
$wsAll.Range("F19").Value = 11146
$wsAll.Range("F21").Value = 5328
